# ABD dip hourly update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 241 - refreshed price/percent/slope values
$ws.Range("D241").Value = 99.8
$ws.Range("E241").Value = 61.67
$ws.Range("F241").Value = 4.61566

# Row 573 - new top entry (MOG-B) replaces previous WDC row
$ws.Range("A573").Value = "MOG-B"
$ws.Range("C573").Value = 195.25
$ws.Range("D573").Value = 325.36
$ws.Range("E573").Value = 66.64
$ws.Range("F573").Value = 9.0632

# Rows 574-584 stay unchanged; rows 585-594 are pushed down from rows 573-584 former
# values (each row's data shifts down by one), and row 595 gets a refreshed date
# and the previous row 594 values.

# Row 585 (previously VALE) becomes WDC (former row 573 data, minus the refreshed slope)
$ws.Range("A585").Value = "WDC"
$ws.Range("C585").Value = 76.97
$ws.Range("D585").Value = 284.1
$ws.Range("E585").Value = 269.1
$ws.Range("F585").Value = 9.58625

# Row 586 (previously SNDK) becomes VALE (former row 585 data)
$ws.Range("A586").Value = "VALE"
$ws.Range("C586").Value = 10.12
$ws.Range("D586").Value = 17.04
$ws.Range("E586").Value = 68.38
$ws.Range("F586").Value = 0.42388

# Row 587 (previously AER) becomes SNDK (former row 586 data)
$ws.Range("A587").Value = "SNDK"
$ws.Range("C587").Value = 46.37
$ws.Range("D587").Value = 630.29
$ws.Range("E587").Value = 1259.26
$ws.Range("F587").Value = 27.57969

# Row 588 (previously LITE) becomes AER (former row 587 data)
$ws.Range("A588").Value = "AER"
$ws.Range("C588").Value = 119.44
$ws.Range("D588").Value = 147.76
$ws.Range("E588").Value = 23.71
$ws.Range("F588").Value = 0.89191

# Row 589 (previously ABVX) becomes LITE (former row 588 data)
$ws.Range("A589").Value = "LITE"
$ws.Range("C589").Value = 119.34
$ws.Range("D589").Value = 583.46
$ws.Range("E589").Value = 388.91
$ws.Range("F589").Value = 54.62861

# Row 590 (previously STX) becomes ABVX (former row 589 data)
$ws.Range("A590").Value = "ABVX"
$ws.Range("C590").Value = 77.28
$ws.Range("D590").Value = 124.54
$ws.Range("E590").Value = 61.15
$ws.Range("F590").Value = 1.53844

# Row 591 (previously AXIA) becomes STX (former row 590 data)
$ws.Range("A591").Value = "STX"
$ws.Range("C591").Value = 159.21
$ws.Range("D591").Value = 431.17
$ws.Range("E591").Value = 170.82
$ws.Range("F591").Value = 11.41482

# Row 592 (previously GLW) becomes AXIA (former row 591 data)
$ws.Range("A592").Value = "AXIA"
$ws.Range("C592").Value = 7.85
$ws.Range("D592").Value = 11.33
$ws.Range("E592").Value = 44.33
$ws.Range("F592").Value = 0.33733

# Row 593 (previously IESC) becomes GLW (former row 592 data)
$ws.Range("A593").Value = "GLW"
$ws.Range("C593").Value = 65.77
$ws.Range("D593").Value = 131.5
$ws.Range("E593").Value = 99.94
$ws.Range("F593").Value = 9.30105

# Row 594 (previously RL) becomes IESC (former row 593 data)
$ws.Range("A594").Value = "IESC"
$ws.Range("C594").Value = 341.11
$ws.Range("D594").Value = 495.49
$ws.Range("E594").Value = 45.26
$ws.Range("F594").Value = 22.4731

# Row 595 (previously MOG-B) becomes RL (former row 594 data) with an updated date.
# The date must stay stored as literal text (matching the rest of column B), so we
# enter it with a leading apostrophe to force text entry and then restore the
# cell's style so no stray number-format is left behind.
$ws.Range("A595").Value = "RL"
$ws.Range("B595").Value = "'2025-08-22"
$ws.Range("B595").Style = "Normal"
$ws.Range("C595").Value = 285.8
$ws.Range("D595").Value = 363.06
$ws.Range("E595").Value = 27.03
$ws.Range("F595").Value = 0.51691
